$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.191.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.519.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.524.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.967.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.136.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.548.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0765"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.803"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0929"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  -1.51%  "
